{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Summary of the change (see commit message \"feat: Add progress #23\n// 'Interview Preparation'\" and the supplied OOXML diff):\n//  1. \"Thank you for reaching out to me.\" -> \"Thank you so much for\n//     reaching out to me.\"\n//  2. \"Not only I have 16 minutes.\" -> \"I have completed booking the\n//     form to further explore the product with you.\"\n//  3. The \"Would you also like to speak with ...\" paragraph is expanded\n//     with \"Downie\" (Jake's last name) and a longer explanation about\n//     the team / COVID-19.\n//  4. A new sign-off block is appended at the end of the letter:\n//       (blank)\n//       (blank)\n//       Thank you,\n//       Hyungmo Gu\n\nconst body = context.document.body;\n\n// Helper: locate the paragraph that currently holds `needle` (an exact,\n// case-sensitive substring of its text) by using Body.search, which is\n// more robust than relying on a fixed paragraph index.\nasync function findParagraph(needle) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find paragraph containing: \" + needle);\n  }\n  const para = results.items[0].paragraphs.getFirst();\n  await context.sync();\n  return para;\n}\n\n// --- 1. \"Thank you for reaching out to me.\" -----------------------------\nconst thanksPara = await findParagraph(\"Thank you for reaching out to me.\");\nthanksPara.getRange(\"Whole\").insertText(\"Thank you\", \"Replace\");\nawait context.sync();\nthanksPara.getRange(\"End\").insertText(\" so much\", \"End\");\nawait context.sync();\nthanksPara.getRange(\"End\").insertText(\" for reaching out to me.\", \"End\");\nawait context.sync();\n\n// --- 2. \"Not only I have 16 minutes.\" ------------------------------------\nconst minutesPara = await findParagraph(\"Not only I have 16 minutes.\");\nminutesPara.getRange(\"Whole\").insertText(\n  \"I have completed booking the form to further explore the product with you.\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 3. \"Would you also like to speak with ...\" --------------------------\nconst inviteParaRuns = [\n  \"Would you also like to speak with Christian Hamm, James Faulkner and Jake\",\n  \" \",\n  \"Downie\",\n  \" at \",\n  \"SiteMax\",\n  \" Systems? They are\",\n  \" amazing people, and they\",\n  \" always\",\n  \" looking to improve their sales and help the company grow, and \",\n  \"during COVID-19 \",\n  \"a help like this \",\n  \"may\",\n  \" be beneficial for the company.\",\n  \"  \",\n];\n\nconst invitePara = await findParagraph(\n  \"Would you also like to speak with Christian Hamm, James Faulkner and Jake\"\n);\ninvitePara.getRange(\"Whole\").insertText(inviteParaRuns[0], \"Replace\");\nawait context.sync();\nfor (let i = 1; i < inviteParaRuns.length; i++) {\n  invitePara.getRange(\"End\").insertText(inviteParaRuns[i], \"End\");\n  await context.sync();\n}\n\n// --- 4. Append the new sign-off block -------------------------------------\nlet anchor = invitePara;\nanchor = anchor.insertParagraph(\"\", \"After\");\nawait context.sync();\nanchor = anchor.insertParagraph(\"\", \"After\");\nawait context.sync();\nanchor = anchor.insertParagraph(\"Thank you,\", \"After\");\nawait context.sync();\nanchor = anchor.insertParagraph(\"Hyungmo\", \"After\");\nawait context.sync();\nanchor.getRange(\"End\").insertText(\" Gu\", \"End\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Summary of the change (see commit message \"feat: Add progress #23\n# 'Interview Preparation'\" and the supplied OOXML diff):\n#  1. \"Thank you for reaching out to me.\" -> \"Thank you so much for\n#     reaching out to me.\"\n#  2. \"Not only I have 16 minutes.\" -> \"I have completed booking the\n#     form to further explore the product with you.\"\n#  3. The \"Would you also like to speak with ...\" paragraph is expanded\n#     with \"Downie\" (Jake's last name) and a longer explanation about\n#     the team / COVID-19.\n#  4. A new sign-off block is appended at the end of the letter:\n#       (blank)\n#       (blank)\n#       Thank you,\n#       Hyungmo Gu\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# --- 1. \"Thank you for reaching out to me.\" -------------------------------\nReplace-Text \"Thank you for reaching out to me.\" \"Thank you so much for reaching out to me.\"\n\n# --- 2. \"Not only I have 16 minutes.\" -------------------------------------\nReplace-Text \"Not only I have 16 minutes.\" \"I have completed booking the form to further explore the product with you.\"\n\n# --- 3. \"Would you also like to speak with ...\" ---------------------------\n$newInviteText = \"Would you also like to speak with Christian Hamm, James Faulkner and Jake Downie at SiteMax Systems? They are amazing people, and they always looking to improve their sales and help the company grow, and during COVID-19 a help like this may be beneficial for the company.  \"\nReplace-Text \"Would you also like to speak with Christian Hamm, James Faulkner and Jake at SiteMax Systems? They are always  \" $newInviteText\n\n# --- 4. Append the new sign-off block --------------------------------------\n#      (blank)\n#      (blank)\n#      Thank you,\n#      Hyungmo Gu\n$endRange = $d.Content\n$endRange.Collapse(0)   # wdCollapseEnd\n$endRange.InsertParagraphAfter()\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n$endRange.Collapse(0)\n\n$count = $d.Paragraphs.Count\n$d.Paragraphs.Item($count - 1).Range.Text = \"Thank you,\"\n$d.Paragraphs.Item($count).Range.Text = \"Hyungmo Gu\"\n"}
